$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Sample_ID value in B2 (100 -> 102)
$ws.Range("B2").Value = 102

# Update the active selection to B2 (was B3)
$ws.Range("B2").Select()
